# ArduinoNode.pptx edit:
#   Duplicate slide 3 (the circuit-board picture that carries the
#   "Rounded Rectangle" component callout) so a fresh copy lands right
#   after it as the new slide 4 - matching the family of per-component
#   callout slides that already follow (old slides 4-7 simply shift down
#   to positions 5-8).
#
#   The original slide 3 is then tidied up: its leftover callout
#   rectangle is removed and its picture is nudged up a touch, leaving a
#   plain "full board" slide ahead of the per-component callouts.

$p = $ppt.ActivePresentation

# 1) Duplicate slide 3; PowerPoint places the copy immediately after the
#    source slide, i.e. it becomes the new slide 4.
$source = $p.Slides.Item(3)
$source.Duplicate() | Out-Null

# 2) Clean up the original slide 3: drop the rounded-rectangle callout
#    and shift the picture's top slightly (450000 EMU -> 436618 EMU).
$original = $p.Slides.Item(3)
$original.Shapes.Item("Rounded Rectangle 5").Delete()
$original.Shapes.Item("Picture 4").Top = 436618 / 914400 * 72
